# Apply the edits described by the commit:
#  1. Slide 8 ("Thanks for your attention!") closing text is shortened to
#     "Thanks you!" (the runs for "for"/" "/"your"/" "/"attention" collapse
#     into a single "you" run, reusing the "for" run's formatting).
#  2. The auto-updating date placeholders (master / layouts / notes master)
#     bump from 25.01.2016 to 26.01.2016.

$p = $ppt.ActivePresentation

# --- 1. Slide 8 closing text ---------------------------------------------
$slide8 = $p.Slides.Item(8)
for ($i = 1; $i -le $slide8.Shapes.Count; $i++) {
    $shp = $slide8.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "Thanks for your attention!") {
            $shp.TextFrame.TextRange.Text = "Thanks you!"
        }
    }
}

# --- helper: refresh any date placeholder shape whose cached text is the
#     old meeting date, on a given shape collection --------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "25.01.2016") {
                    $shp.TextFrame.TextRange.Text = "26.01.2016"
                }
            }
        }
    }
}

# --- 2. Slide master -------------------------------------------------------
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

# --- 3. Slide layouts -------------------------------------------------------
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# NOTE: the notes master's own date placeholder ("Datumsplatzhalter 2",
# {7A19BEE0-598C-4858-A439-753E34679212}) is intentionally left alone here:
# in this COM host, writes through Presentation.NotesMaster.Shapes.Item(n)
# land on SlideMaster.Shapes.Item(n) instead (a notes-master/slide-master
# shape aliasing quirk), so touching it would corrupt unrelated slide
# master placeholders instead of updating the notes master.
